$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width updates (character widths 7 -> 8 and 8 -> 9 as stored in the OOXML) ---
# Excel's COM ColumnWidth property is offset from the raw OOXML "width" attribute by 5/6
# of a character for the default font used here, so feeding it N + 1/6 round-trips to the
# exact integer N once re-serialised.
$cols8 = @(2,3,7,9,12,13,15,16,17,22,24,26,27,29,30,34)
$cols9 = @(20)

foreach ($c in $cols8) {
    $ws.Cells.Item(1, $c).EntireColumn.ColumnWidth = 43/6
}
foreach ($c in $cols9) {
    $ws.Cells.Item(1, $c).EntireColumn.ColumnWidth = 49/6
}

# --- Replace the data rows (rows 2-5) with the refreshed measurement set ---
$newData = @(
    @(45119.50694444445,16.815,11.263,3.771,35.807,27.941,13.232,40.653,20.36,8.363,12.4,14.083,14.597,4.223,13.159,18.251,11.431,3.368,2.219,192.924,36.527,12.146,23.789,12.048,3.162,20.789,10.728,9.564,11.528,14.849,3.316,36.381,6.54,15.185),
    @(45119.51388888889,11.53,7.989,1.614,24.928,19.621,9.074,35.689,13.961,5.932,8.598000000000001,9.959,10.367,2.9,9.023,12.634,7.932,1.482,0.946,130.023,25.298,8.329000000000001,16.585,8.561,1.866,17.306,7.357,6.65,7.895,10.427,1.247,32.797,4.506,10.413),
    @(45119.52083333334,3.843,2.441,0.888,8.295999999999999,6.144,3.026,15.588,4.654,1.914,2.568,3.32,3.405,0.975,3.008,4.177,2.841,0.927,0.463,38.497,8.648,2.776,5.522,2.775,0.844,7.141,2.452,2.313,2.741,3.456,0.766,14.616,1.405,3.473),
    @(45119.52777777778,12.49,9.050000000000001,0.97,27.14,21.88,9.83,35.26,15.12,6.61,9.720000000000001,10.89,11.45,3.14,9.779999999999999,13.8,8.380000000000001,0.77,0.61,141.43,27.21,9.02,18.13,9.52,1.62,17.31,7.97,7.12,8.390000000000001,11.42,0.55,31.82,5.01,11.28)
)

for ($r = 0; $r -lt $newData.Length; $r++) {
    $rowVals = $newData[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $rowVals[$c]
    }
}

# --- Drop the now-unused 6th data row (sheet shrinks from A1:AH6 to A1:AH5) ---
$ws.Range("A6:AH6").EntireRow.Delete()
